# Add the final "Viva exam / Viva and public defense" milestone row to the
# "Chronology of work" sheet (this was previously a stub row with only an
# empty, pre-formatted B20 cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chronology of work")
$ws.Activate()

# A20: date of the viva (04-Sep-2023), formatted like the other date cells
# in column A (built-in date format "d-mmm-yy").
$ws.Range("A20").Value = "9/4/2023"
$ws.Range("A20").NumberFormat = "d-mmm-yy"

# B20 already carries the "green milestone text" style (s="9"); just fill
# in the text, keeping that formatting.
$ws.Range("B20").Value = "Viva exam"

# C20 (Type) and D20 (Comments) for the new row.
$ws.Range("C20").Value = "Viva exam"
$ws.Range("D20").Value = "Viva and public defense"

# Move the selection down to A21, as if the user had just finished typing
# this row and pressed Enter.
$ws.Range("A21").Select() | Out-Null
